$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# N2: update report date text
$ws.Range("N2").Value = "2020-09-30 00:00:00"

# Numeric columns O2:AC2 - update values
$ws.Range("O2").Value = 44467642.2
$ws.Range("P2").Value = 122.2918776159
$ws.Range("Q2").Value = 172447959.65
$ws.Range("R2").Value = 474.254620512
$ws.Range("S2").Value = 47160015.48
$ws.Range("T2").Value = 129.6962590349
$ws.Range("U2").Value = -849484509.8099999
$ws.Range("V2").Value = -2336.1943779934
$ws.Range("W2").Value = 4491325.39
$ws.Range("X2").Value = 12.351736853
$ws.Range("Y2").Value = 127400835.2
$ws.Range("Z2").Value = 350.3690903233
$ws.Range("AA2").Value = 841468742.47
$ws.Range("AB2").Value = 2314.1499611986
$ws.Range("AC2").Value = 36361893.42

# AD2 was an empty inline string cell; now becomes a numeric value
$ws.Range("AD2").Value = 1092.0130062434
